$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "84×50="
$t.Cell(1,2).Range.Text = "78×67="
$t.Cell(1,3).Range.Text = "15×24="
$t.Cell(1,4).Range.Text = "22×36="
$t.Cell(1,5).Range.Text = "69×67="
$t.Cell(5,1).Range.Text = "98×65="
$t.Cell(5,2).Range.Text = "99×18="
$t.Cell(5,3).Range.Text = "89×11="
$t.Cell(5,4).Range.Text = "42×55="
$t.Cell(5,5).Range.Text = "42×82="
$t.Cell(10,1).Range.Text = "31×71="
$t.Cell(10,2).Range.Text = "21×47="
$t.Cell(10,3).Range.Text = "70×32="
$t.Cell(10,4).Range.Text = "46×43="
$t.Cell(10,5).Range.Text = "56×36="
$t.Cell(15,1).Range.Text = "34×53="
$t.Cell(15,2).Range.Text = "58×17="
$t.Cell(15,3).Range.Text = "50×96="
$t.Cell(15,4).Range.Text = "74×40="
$t.Cell(15,5).Range.Text = "28×18="
$t.Cell(20,1).Range.Text = "22×19="
$t.Cell(20,2).Range.Text = "95×44="
$t.Cell(20,3).Range.Text = "57×14="
$t.Cell(20,4).Range.Text = "48×23="
$t.Cell(20,5).Range.Text = "73×37="
